$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new BOM row above the MOSFET (Q1) row, pushing the rows below it
# down by one, and give the new row the same formatting as the rest of the
# BOM table (copy it from the row that used to be row 42, now row 43).
$ws.Rows.Item(42).Insert()
$ws.Range("A43:F43").Copy()
$ws.Range("A42:F42").PasteSpecial(-4122)

# Fill in the new TVS diode line item (D4 / CDSOD323-T03SC / Bourns).
# Description, part number, manufacturer and URL are entered before the
# designator so the shared-string table picks up the same ordering as the
# authored workbook.
$ws.Range("B42").Value = "TVS DIODE 3,3V 10,9V SOD323"
$ws.Range("C42").Value = "CDSOD323-T03SC"
$ws.Range("D42").Value = "BOURNS INC"
$ws.Range("E42").Value = "https://octopart.com/cdsod323-t03sc-bourns-10487153?r=sp"
$ws.Range("A42").Value = "D4"
$ws.Range("F42").Value = 1

# Scroll/select the same way the author left the sheet.
$ws.Range("A29").Select()
$ws.Range("B39").Select()
